$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted above the existing row 169,
# shifting the former rows 169-172 down to 170-173 (data otherwise intact).
$ws.Rows.Item(169).Insert()

# Populate the newly inserted row 169 with the new observation.
$ws.Cells.Item(169, 1).Value = 10
$ws.Cells.Item(169, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(169, 3).Value = "La Araucanía"
$ws.Cells.Item(169, 4).Value = 44448
$ws.Cells.Item(169, 5).Value = 9
$ws.Cells.Item(169, 6).Value = 100112001
$ws.Cells.Item(169, 7).Value = "Berenjena"
$ws.Cells.Item(169, 8).Value = "Sin especificar"
$ws.Cells.Item(169, 9).Value = "Primera"
$ws.Cells.Item(169, 10).Value = 75
$ws.Cells.Item(169, 11).Value = 12000
$ws.Cells.Item(169, 12).Value = 13000
$ws.Cells.Item(169, 13).Value = 12533
$ws.Cells.Item(169, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(169, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(169, 16).Value = 209
$ws.Cells.Item(169, 17).Value = 60
$ws.Cells.Item(169, 18).Value = "Hortaliza"
